# Update factsheets with text edits from COMM
#
# Converts numeric "No. of 990 Filers w/ Gov Grants" cells (column B on the
# detail sheets, column A on the Overall sheet) from real numbers into text
# values (so they keep the "1,446"-style thousands separator as literal
# text), fixes the "Greenlee County" placeholder row on the County sheet,
# and appends a new "Total" row to the County sheet.

function Set-TextValue($ws, $addr, $val) {
    # Forcing NumberFormat to Text ("@") before assigning the value stops
    # Excel's autodetection from turning numeric-looking strings (e.g.
    # "1,446", "12", "0.00%", "$0") back into real numbers. ClearFormats()
    # afterwards drops the temporary Text number format again so the cell
    # ends up without any explicit style, matching the original formatting.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overall sheet: A2 1446 -> "1,446"
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall "A2" "1,446"

# ---------------------------------------------------------------------
# County sheet
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")

Set-TextValue $wsCounty "B2" "12"
Set-TextValue $wsCounty "B3" "26"
Set-TextValue $wsCounty "B4" "64"
Set-TextValue $wsCounty "B5" "12"
Set-TextValue $wsCounty "B6" "10"
Set-TextValue $wsCounty "B7" "4"
Set-TextValue $wsCounty "B8" "847"
Set-TextValue $wsCounty "B9" "26"
Set-TextValue $wsCounty "B10" "20"
Set-TextValue $wsCounty "B11" "267"
Set-TextValue $wsCounty "B12" "35"
Set-TextValue $wsCounty "B13" "15"
Set-TextValue $wsCounty "B14" "83"
Set-TextValue $wsCounty "B15" "25"

# Greenlee County row (row 16) placeholder values updated
Set-TextValue $wsCounty "B16" "0.00%"
Set-TextValue $wsCounty "C16" "$0"
Set-TextValue $wsCounty "D16" "0.00%"
Set-TextValue $wsCounty "E16" "0.00%"
Set-TextValue $wsCounty "F16" "0.00%"

# New Total row (row 17)
Set-TextValue $wsCounty "A17" "Total"
Set-TextValue $wsCounty "B17" "1,446"
Set-TextValue $wsCounty "C17" "$3,710,995,891"
Set-TextValue $wsCounty "D17" "8.70%"
Set-TextValue $wsCounty "E17" "-12.71%"
Set-TextValue $wsCounty "F17" "65.98%"

# ---------------------------------------------------------------------
# Congressional District sheet
# ---------------------------------------------------------------------
$wsCD = $wb.Worksheets.Item("Congressional District")

Set-TextValue $wsCD "B2" "192"
Set-TextValue $wsCD "B3" "202"
Set-TextValue $wsCD "B4" "284"
Set-TextValue $wsCD "B5" "170"
Set-TextValue $wsCD "B6" "79"
Set-TextValue $wsCD "B7" "145"
Set-TextValue $wsCD "B8" "218"
Set-TextValue $wsCD "B9" "92"
Set-TextValue $wsCD "B10" "64"
Set-TextValue $wsCD "B11" "1,446"

# ---------------------------------------------------------------------
# Size sheet
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")

Set-TextValue $wsSize "B2" "349"
Set-TextValue $wsSize "B3" "453"
Set-TextValue $wsSize "B4" "245"
Set-TextValue $wsSize "B5" "166"
Set-TextValue $wsSize "B6" "191"
Set-TextValue $wsSize "B7" "42"
Set-TextValue $wsSize "B8" "1,446"

# ---------------------------------------------------------------------
# Subsector sheet
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")

Set-TextValue $wsSub "B2" "124"
Set-TextValue $wsSub "B3" "352"
Set-TextValue $wsSub "B4" "62"
Set-TextValue $wsSub "B5" "123"
Set-TextValue $wsSub "B6" "10"
Set-TextValue $wsSub "B7" "378"
Set-TextValue $wsSub "B8" "11"
Set-TextValue $wsSub "B9" "98"
Set-TextValue $wsSub "B10" "33"
Set-TextValue $wsSub "B11" "249"
Set-TextValue $wsSub "B12" "6"
Set-TextValue $wsSub "B13" "1,446"
